# Generate Report for Handoff
# Adds two new source files (6cece956-27ec-46ad-a6b1-780773bbaed1.md and
# a561822c-30dd-4e4a-9619-cbc7760525eb.md) as new rows 4 & 5 across the
# Overview, zh-cn and de-de sheets, mirroring the existing "Ready for
# handoff" row pattern (row 3 / af120165.. style row).
#
# NOTE: all text assignments use a leading "'" so that values which look
# like booleans ("True"/"False") or dates/numbers stay plain text (t="s"),
# matching how the existing workbook stores them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Append two rows by copying the last existing data row (row 3) so that
# cell styles (s="0"/"1"/"2") are preserved exactly.
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(5).Insert()

$ws.Range("A4").Value = "'6cece956-27ec-46ad-a6b1-780773bbaed1.md"
$ws.Range("B4").Value = "'e2e\6cece956-27ec-46ad-a6b1-780773bbaed1.md"
$ws.Range("C4").Value = "'.md"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "'Ready for handoff"
$ws.Range("F4").Value = "'Ready for handoff"
$ws.Range("G4").Value = "'2016-08-12 16:50:22"

$ws.Range("A5").Value = "'a561822c-30dd-4e4a-9619-cbc7760525eb.md"
$ws.Range("B5").Value = "'e2e\a561822c-30dd-4e4a-9619-cbc7760525eb.md"
$ws.Range("C5").Value = "'.md"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "'Ready for handoff"
$ws.Range("F5").Value = "'Ready for handoff"
$ws.Range("G5").Value = "'2016-08-12 16:50:22"

$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/master/e2e/6cece956-27ec-46ad-a6b1-780773bbaed1.md", "", "", "e2e\6cece956-27ec-46ad-a6b1-780773bbaed1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/master/e2e/a561822c-30dd-4e4a-9619-cbc7760525eb.md", "", "", "e2e\a561822c-30dd-4e4a-9619-cbc7760525eb.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(5).Insert()

$ws.Range("A4").Value = "'6cece956-27ec-46ad-a6b1-780773bbaed1.md"
$ws.Range("B4").Value = "'.md"
$ws.Range("C4").Value = "'Ready for handoff"
$ws.Range("D4").Value = "'e2e"
$ws.Range("E4").Value = "'ht"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "'6cece956-27ec-46ad-a6b1-780773bbaed1.1852d253f225ccc10a0c4dae60242cea4b9bb265.zh-cn.xlf"
$ws.Range("H4").Value = "'2016-08-12 16:50:16"
$ws.Range("I4").Value = "'"
$ws.Range("J4").Value = "'"
$ws.Range("K4").Value = "'0001-01-01 00:00:00"
$ws.Range("L4").Value = "'"
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = "'"
$ws.Range("O4").Value = "'False"
$ws.Range("P4").Value = "'"

$ws.Range("A5").Value = "'a561822c-30dd-4e4a-9619-cbc7760525eb.md"
$ws.Range("B5").Value = "'.md"
$ws.Range("C5").Value = "'Ready for handoff"
$ws.Range("D5").Value = "'e2e"
$ws.Range("E5").Value = "'ht"
$ws.Range("F5").Value = "'False"
$ws.Range("G5").Value = "'a561822c-30dd-4e4a-9619-cbc7760525eb.4498ec71c2b42802aae9dc7d8bb33baf9cb6e553.zh-cn.xlf"
$ws.Range("H5").Value = "'2016-08-12 16:50:16"
$ws.Range("I5").Value = "'"
$ws.Range("J5").Value = "'"
$ws.Range("K5").Value = "'0001-01-01 00:00:00"
$ws.Range("L5").Value = "'"
$ws.Range("M5").Value = "'True"
$ws.Range("N5").Value = "'"
$ws.Range("O5").Value = "'False"
$ws.Range("P5").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/6cece956-27ec-46ad-a6b1-780773bbaed1.md", "", "", "6cece956-27ec-46ad-a6b1-780773bbaed1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/a561822c-30dd-4e4a-9619-cbc7760525eb.md", "", "", "a561822c-30dd-4e4a-9619-cbc7760525eb.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Copy()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(5).Insert()

$ws.Range("A4").Value = "'6cece956-27ec-46ad-a6b1-780773bbaed1.md"
$ws.Range("B4").Value = "'.md"
$ws.Range("C4").Value = "'Ready for handoff"
$ws.Range("D4").Value = "'e2e"
$ws.Range("E4").Value = "'ht"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "'6cece956-27ec-46ad-a6b1-780773bbaed1.1852d253f225ccc10a0c4dae60242cea4b9bb265.de-de.xlf"
$ws.Range("H4").Value = "'2016-08-12 16:50:22"
$ws.Range("I4").Value = "'"
$ws.Range("J4").Value = "'"
$ws.Range("K4").Value = "'0001-01-01 00:00:00"
$ws.Range("L4").Value = "'"
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = "'"
$ws.Range("O4").Value = "'False"
$ws.Range("P4").Value = "'"

$ws.Range("A5").Value = "'a561822c-30dd-4e4a-9619-cbc7760525eb.md"
$ws.Range("B5").Value = "'.md"
$ws.Range("C5").Value = "'Ready for handoff"
$ws.Range("D5").Value = "'e2e"
$ws.Range("E5").Value = "'ht"
$ws.Range("F5").Value = "'False"
$ws.Range("G5").Value = "'a561822c-30dd-4e4a-9619-cbc7760525eb.4498ec71c2b42802aae9dc7d8bb33baf9cb6e553.de-de.xlf"
$ws.Range("H5").Value = "'2016-08-12 16:50:22"
$ws.Range("I5").Value = "'"
$ws.Range("J5").Value = "'"
$ws.Range("K5").Value = "'0001-01-01 00:00:00"
$ws.Range("L5").Value = "'"
$ws.Range("M5").Value = "'True"
$ws.Range("N5").Value = "'"
$ws.Range("O5").Value = "'False"
$ws.Range("P5").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/6cece956-27ec-46ad-a6b1-780773bbaed1.md", "", "", "6cece956-27ec-46ad-a6b1-780773bbaed1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/a561822c-30dd-4e4a-9619-cbc7760525eb.md", "", "", "a561822c-30dd-4e4a-9619-cbc7760525eb.md") | Out-Null
